# Lab 1 Report.docx edit:
#  1. Remove the direct w:sz/w:szCs (28 half-points => 14pt) run/paragraph
#     formatting from the "To be familiar with using digital multimeter
#     and oscilloscope" introduction bullet (paragraph 21, w14:paraId
#     "7FA64C73"), letting it fall back to the style/doc defaults.
#  2. Remove the stray <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> from
#     the pPr of the lone page-break paragraph (paragraph 327, w14:paraId
#     "1FDEE954") right after "Figure 7. Measure current versus voltage
#     for the diode".
#
# Both paragraphs are rewritten in-place via Range.InsertXML, reproducing
# every attribute (paraId, textId, rsids, run rsids, etc.) exactly so that
# only the targeted child elements disappear.

$d = $word.ActiveDocument

# --- 1. "To be familiar with using digital multimeter and oscilloscope" ---
$introPara = $d.Paragraphs(21)
$introXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
          '<w:body>' +
            '<w:p w14:paraId="7FA64C73" w14:textId="34720C3B" w:rsidR="008C3E65" w:rsidRDefault="008C3E65" w:rsidP="008C3E65">' +
              '<w:pPr>' +
                '<w:pStyle w:val="a3"/>' +
                '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>' +
                '<w:ind w:leftChars="0"/>' +
              '</w:pPr>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>T</w:t></w:r>' +
              '<w:r><w:t>o be familiar with using digital multimeter</w:t></w:r>' +
              '<w:r w:rsidR="0078607F"><w:t xml:space="preserve"> and </w:t></w:r>' +
              '<w:r w:rsidR="0078607F"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>o</w:t></w:r>' +
              '<w:r w:rsidR="0078607F"><w:t>scilloscope</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$introPara.Range.InsertXML($introXml)

# --- 2. page-break paragraph after Figure 7 ---
$pageBreakPara = $d.Paragraphs(327)
$pbXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
          '<w:body>' +
            '<w:p w14:paraId="1FDEE954" w14:textId="6F5D77C8" w:rsidR="00C52D7A" w:rsidRDefault="00776C35" w:rsidP="00776C35">' +
              '<w:pPr>' +
                '<w:widowControl/>' +
              '</w:pPr>' +
              '<w:r><w:br w:type="page"/></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$pageBreakPara.Range.InsertXML($pbXml)
